$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "Sheet1" -> "Kane Williamson"
$ws.Name = "Kane Williamson"

# 2. Insert a new column before column A (shifts teamName..result from A:L to B:M)
$ws.Columns.Item(1).Insert()

# Helper: write a value as TEXT even when it looks numeric (e.g. "66", "84.21"),
# mirroring the source file where every cell is stored as t="str". Plain
# Range.Value assignment auto-converts numeric-looking strings into real
# numbers, so instead we push the literal string in via a formula and then
# collapse the formula down to its cached value with a values-only paste.
function Set-TextValue($cell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# 3. New header cell
$ws.Range("A1").Value = "matchNo"

# 4. matchNo for the pre-existing row (now row 2)
Set-TextValue $ws.Range("A2") "23rd"

# 5. New row 3 (20th match)
Set-TextValue $ws.Range("A3") "20th"
$ws.Range("B3").Value = "Sunrisers Hyderabad"
$ws.Range("C3").Value = "Kane Williamson"
Set-TextValue $ws.Range("E3") "66"
Set-TextValue $ws.Range("F3") "51"
Set-TextValue $ws.Range("G3") "8"
Set-TextValue $ws.Range("H3") "0"
Set-TextValue $ws.Range("I3") "129.41"
$ws.Range("J3").Value = "Delhi Capitals"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 25"
$ws.Range("M3").Value = "Match tied (Capitals won the one-over eliminator)"

# 6. New row 4 (14th match)
Set-TextValue $ws.Range("A4") "14th"
$ws.Range("B4").Value = "Sunrisers Hyderabad"
$ws.Range("C4").Value = "Kane Williamson"
Set-TextValue $ws.Range("E4") "16"
Set-TextValue $ws.Range("F4") "19"
Set-TextValue $ws.Range("G4") "0"
Set-TextValue $ws.Range("H4") "0"
Set-TextValue $ws.Range("I4") "84.21"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Chennai"
$ws.Range("L4").Value = "April 21"
$ws.Range("M4").Value = "Sunrisers won by 9 wickets (with 8 balls remaining)"
